$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 20 (M6S, count 1) entirely; rows below shift up.
$ws.Rows.Item(20).Delete()
